$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700.0
# half-an-EMU nudge so float32-precision COM properties (Left/Top/Width/Height
# are exposed as single-precision floats) round-trip to the exact EMU integer
# PowerPoint would have written instead of landing 1 EMU short.
$epsilon = 0.5 / $emuPerPt

function ToPt([double]$emu) {
    return ($emu / $emuPerPt) + $epsilon
}

# ---------------------------------------------------------------------------
# 1) "CAP 1" textbox (shape id=5) becomes the "Cap1.doc.docx" hyperlink box.
# ---------------------------------------------------------------------------
$capShape = $s.Shapes.Item(2)

$capRange = $capShape.TextFrame.TextRange
$capRange.Text = "Cap1.doc.docx"
$capLink = $capRange.ActionSettings.Item(1).Hyperlink
$capLink.Address = "Cap1.doc.docx"

$capShape.Width = ToPt 1820256
$capShape.Height = ToPt 1200329
$capShape.Left = ToPt 9802026
$capShape.Top = ToPt 449938

# ---------------------------------------------------------------------------
# 2) New "CaixaDeTexto 2" textbox (must land on shape id=3) linking to
#    "Cap2.doc.docx". Shape ids are handed out as the lowest free integer, so
#    we burn id=2 on a throwaway duplicate (immediately deleted) before
#    duplicating the (already-updated) CAP1 box for real.
# ---------------------------------------------------------------------------
$throwaway = $s.Shapes.Item(1).Duplicate()
$throwawayShape = $throwaway.Item(1)
$throwawayShape.Delete()

$newRange = $capShape.Duplicate()
$cap2Shape = $newRange.Item(1)
$cap2Shape.Name = "CaixaDeTexto 2"

$cap2Range = $cap2Shape.TextFrame.TextRange
$cap2Range.Text = "Cap2.doc.docx"
$cap2Range.Font.Size = 32
$cap2Link = $cap2Range.ActionSettings.Item(1).Hyperlink
$cap2Link.Address = "Cap2.doc.docx"

$cap2Shape.Width = ToPt 1580972
$cap2Shape.Height = ToPt 1077218
$cap2Shape.Left = ToPt 9827664
$cap2Shape.Top = ToPt 2196269
